# Update header text on the two existing sheets
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the existing sheets
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy header cell style (bold, centered, bordered) from an existing header cell
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Data rows
$rows = @(
    @(45459.99999999999, 73,  0.7836942644172917, 144.1213094654169),
    @(45473.99999999999, 79,  8.257630732060889,  145.597359776045),
    @(45494.99999999999, 87,  14.71237565402322,  159.1192966232402),
    @(45501.99999999999, 90,  16.16303684707027,  158.7868230402784),
    @(45508.99999999999, 93,  25.49652951417849,  164.8681259927558),
    @(45515.99999999999, 96,  29.39641799523318,  165.8481567648835),
    @(45522.99999999999, 99,  32.12215351016929,  167.7655894081431),
    @(45529.99999999999, 102, 30.70333865230726,  173.8124616708128),
    @(45536.99999999999, 105, 31.72957933634249,  173.6912373767797),
    @(45543.99999999999, 108, 35.35553355935718,  179.0710661872735),
    @(45550.99999999999, 110, 35.47795252921959,  180.5407455365333),
    @(45557.99999999999, 113, 38.55359564287714,  189.4048935191765)
)

$r = 2
foreach ($row in $rows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Copy date style (numFmt) from existing date column onto the new ds column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A13").PasteSpecial(-4122)

Write-Host "Edit complete"
